$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..20: Team, B(calc_PPDA), C(calc_OPPDA), D(calc_FieldTilt_Pct), E(calc_HighError_Rate), F(calc_Directness), G(calc_BigChance_Diff)
$data = @(
    @("Barcelona", 3.889784946236559, 8.716510903426791, 0.7921889568009945, 12, 0.04400647948164147, 55),
    @("Real Madrid", 5.164348925410872, 7.045523520485585, 0.6814670017685935, 12, 0.05884992148358633, 38),
    @("Atlético Madrid", 4.93520782396088, 8.992818671454218, 0.5458039052762775, 15, 0.06744302949061662, 29),
    @("Villarreal", 5.572571428571429, 7.23869801084991, 0.405093926472283, 29, 0.08306567726215527, 23),
    @("Real Betis", 5.656470588235294, 8.105924596050269, 0.5216118074264964, 12, 0.07946646627841443, 9),
    @("Espanyol", 5.88496583143508, 6.791907514450867, 0.4016571583808747, 10, 0.09650195694716243, 6),
    @("Celta Vigo", 5.711168164313222, 9.212612612612613, 0.4927824267782427, 18, 0.07421977445580907, 0),
    @("Real Sociedad", 4.335758835758836, 7.474074074074074, 0.480319803198032, 9, 0.0735734186358477, 11),
    @("Osasuna", 5.334928229665072, 6.573012939001848, 0.421151674821227, 11, 0.09909289241015042, -11),
    @("Athletic Club", 6.029274004683841, 5.373292867981791, 0.5984109283523836, 20, 0.07825727854635557, 0),
    @("Getafe", 4.559162303664921, 4.326, 0.4154395359214637, 13, 0.1313390313390313, -17),
    @("Girona FC", 5.493638676844784, 8.45362563237774, 0.4401503582755785, 20, 0.07318728189220629, -14),
    @("Sevilla", 3.539915966386554, 7.943661971830986, 0.5478051713770294, 19, 0.08658314811072944, -33),
    @("Deportivo Alavés", 4.206217616580311, 5.849779086892489, 0.4948209160114813, 12, 0.07809485122372296, -2),
    @("Elche", 3.719399538106236, 9.204991087344029, 0.5947818055876241, 17, 0.05305506663314056, -34),
    @("Mallorca", 5.482535575679172, 7.007751937984496, 0.3620145631067961, 11, 0.1101431980906921, -11),
    @("Valencia", 5.064070351758794, 7.863805970149254, 0.433969185619956, 12, 0.0886342130165072, -10),
    @("Rayo Vallecano", 4.011764705882353, 9.272536687631026, 0.537260151717983, 18, 0.08896244957070445, -1),
    @("Levante UD", 5.365644171779141, 6.041366906474821, 0.344306049822064, 11, 0.1054054054054054, -9)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $ws.Cells.Item($row, 6).Value = $item[5]
    $ws.Cells.Item($row, 7).Value = $item[6]
    $row++
}
